$d = $word.ActiveDocument

# Locate the "US Military Networks" paragraph (the last of the three
# bullet-less "Targets:" entries).
$findRange = $d.Content
$findRange.Find.ClearFormatting()
$findRange.Find.Execute("US Military Networks", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$targetStart = $findRange.Start

$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Start -eq $targetStart) {
        $targetIndex = $i
        break
    }
}

# Range covering just the paragraph's text, excluding its trailing
# paragraph mark.
$rng = $d.Paragraphs.Item($targetIndex).Range
[void]$rng.MoveEnd(1, -1)
$rng.Collapse(0)

# Temporarily insert a marker character so we have a valid (non zero-length)
# range to anchor a bookmark on - this avoids an edge case where adding a
# bookmark directly on a zero-length range sitting exactly at a paragraph
# boundary snaps to the wrong location.
$rng.InsertAfter("X")

# Relocate the document's "_GoBack" bookmark (Word's "last edit location"
# marker) to sit right after "US Military Networks", matching where Word
# would leave it after the new text below is typed in.
$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()
$d.Bookmarks.Add("_GoBack", $rng)

# Remove the temporary marker character, collapsing the bookmark back down
# to a zero-length range at the correct position.
$rng.Text = ""

# Insert the new "Targets:" entry as its own paragraph, right after
# "US Military Networks", matching the bold formatting used by the other
# entries in that list.
$insertionPoint = $d.Range($rng.End, $rng.End)
$insertionPoint.InsertParagraphAfter()

$newParagraph = $d.Paragraphs.Item($targetIndex + 1)
$newParagraph.Range.Text = "Japanese, Vietnamese, Iran, Filipino, and Tibetan sites and networks"
